$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 changes from the text "R40" to the text "1". A plain
# Range("B11").Value = "1" assignment would be auto-coerced to the
# *number* 1 by Excel's type inference (since "1" parses as numeric),
# which would lose the original text ("t=""s""") cell type and also the
# cell's existing style if we forced text via NumberFormat="@". To keep
# the value as literal text "1" while leaving B11's style untouched, we
# stage the text in a scratch cell (as a formula that evaluates to the
# string "1", so the scratch cell's own formatting/style never changes),
# then copy just the resulting value onto B11 with PasteSpecial, and
# finally clean up the scratch cell.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=""1"""
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
